$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("general")
$ws2 = $wb.Worksheets.Item("all")

# --- Sheet "general": update header + values ---
$ws1.Range("A1").Value = "funct"
$ws1.Range("B1").Value = "ectrl_area"
$ws1.Range("C1").Value = "ses_area"
$ws1.Range("A2").Value = "ATCOs in OPS"
$ws1.Range("B2").Value = 149000
$ws1.Range("C2").Value = 174600
$ws1.Range("A3").Value = "Support Staff"
$ws1.Range("B3").Value = 78600
$ws1.Range("C3").Value = 96700
$ws1.Range("A4").Value = "Average all staff"
$ws1.Range("B4").Value = 101300
$ws1.Range("C4").Value = 122900

# --- Sheet "all": update header + values ---
$ws2.Range("A1").Value = "ansp"
$ws2.Range("B1").Value = "atco_ops"
$ws2.Range("C1").Value = "support"
$ws2.Range("D1").Value = "all"
$ws2.Range("A2").Value = "Albcontrol"
$ws2.Range("B2").Value = 25157
$ws2.Range("C2").Value = 12004
$ws2.Range("D2").Value = 14278
$ws2.Range("A3").Value = "ANS CR"
$ws2.Range("B3").Value = 92480
$ws2.Range("C3").Value = 48890
$ws2.Range("D3").Value = 59632
$ws2.Range("A4").Value = "ARMATS"
$ws2.Range("B4").Value = 26608
$ws2.Range("C4").Value = 15573
$ws2.Range("D4").Value = 18564
$ws2.Range("A5").Value = "Austro Control"
$ws2.Range("B5").Value = 212241
$ws2.Range("C5").Value = 144827
$ws2.Range("D5").Value = 168198
$ws2.Range("A6").Value = "Avinor (Continental)"
$ws2.Range("B6").Value = 200980
$ws2.Range("C6").Value = 100844
$ws2.Range("D6").Value = 141867
$ws2.Range("A7").Value = "BHANSA"
$ws2.Range("B7").Value = 33229
$ws2.Range("C7").Value = 27790
$ws2.Range("D7").Value = 29339
$ws2.Range("A8").Value = "BULATSA"
$ws2.Range("B8").Value = 97975
$ws2.Range("C8").Value = 43886
$ws2.Range("D8").Value = 57983
$ws2.Range("A9").Value = "Croatia Control"
$ws2.Range("B9").Value = 114253
$ws2.Range("C9").Value = 55490
$ws2.Range("D9").Value = 75575
$ws2.Range("A10").Value = "DCAC Cyprus"
$ws2.Range("B10").Value = 103333
$ws2.Range("C10").Value = 65184
$ws2.Range("D10").Value = 82867
$ws2.Range("A11").Value = "DFS"
$ws2.Range("B11").Value = 259922
$ws2.Range("C11").Value = 108590
$ws2.Range("D11").Value = 154303
$ws2.Range("A12").Value = "DHMI"
$ws2.Range("B12").Value = 48531
$ws2.Range("C12").Value = 13627
$ws2.Range("D12").Value = 21911
$ws2.Range("A13").Value = "DSNA"
$ws2.Range("B13").Value = 135757
$ws2.Range("C13").Value = 101783
$ws2.Range("D13").Value = 114204
$ws2.Range("A14").Value = "EANS"
$ws2.Range("B14").Value = 118386
$ws2.Range("C14").Value = 44161
$ws2.Range("D14").Value = 68337
$ws2.Range("A15").Value = "ENAIRE"
$ws2.Range("B15").Value = 189385
$ws2.Range("C15").Value = 100512
$ws2.Range("D15").Value = 136938
$ws2.Range("A16").Value = "ENAV"
$ws2.Range("B16").Value = 160131
$ws2.Range("C16").Value = 107708
$ws2.Range("D16").Value = 131787
$ws2.Range("A17").Value = "Fintraffic ANS"
$ws2.Range("B17").Value = 121324
$ws2.Range("C17").Value = 91881
$ws2.Range("D17").Value = 106711
$ws2.Range("A18").Value = "HASP"
$ws2.Range("B18").Value = 79373
$ws2.Range("C18").Value = 56721
$ws2.Range("D18").Value = 64127
$ws2.Range("A19").Value = "HungaroControl"
$ws2.Range("B19").Value = 117751
$ws2.Range("C19").Value = 45005
$ws2.Range("D19").Value = 62468
$ws2.Range("A20").Value = "IAA"
$ws2.Range("B20").Value = 136394
$ws2.Range("C20").Value = 110856
$ws2.Range("D20").Value = 123423
$ws2.Range("A21").Value = "LFV"
$ws2.Range("B21").Value = 208646
$ws2.Range("C21").Value = 94208
$ws2.Range("D21").Value = 147432
$ws2.Range("A22").Value = "LGS"
$ws2.Range("B22").Value = 68079
$ws2.Range("C22").Value = 33846
$ws2.Range("D22").Value = 42079
$ws2.Range("A23").Value = "LPS"
$ws2.Range("B23").Value = 111971
$ws2.Range("C23").Value = 36682
$ws2.Range("D23").Value = 53725
$ws2.Range("A24").Value = "LVNL"
$ws2.Range("B24").Value = 169343
$ws2.Range("C24").Value = 130050
$ws2.Range("D24").Value = 137302
$ws2.Range("A25").Value = "MATS"
$ws2.Range("B25").Value = 95426
$ws2.Range("C25").Value = 58823
$ws2.Range("D25").Value = 70659
$ws2.Range("A26").Value = "M-NAV"
$ws2.Range("B26").Value = 62311
$ws2.Range("C26").Value = 26413
$ws2.Range("D26").Value = 34070
$ws2.Range("A27").Value = "MOLDATSA"
$ws2.Range("B27").Value = 29717
$ws2.Range("C27").Value = 16586
$ws2.Range("D27").Value = 20429
$ws2.Range("A28").Value = "MUAC"
$ws2.Range("B28").Value = 342107
$ws2.Range("C28").Value = 201643
$ws2.Range("D28").Value = 248610
$ws2.Range("A29").Value = "NATS (Continental)"
$ws2.Range("B29").Value = 188931
$ws2.Range("C29").Value = 79441
$ws2.Range("D29").Value = 115172
$ws2.Range("A30").Value = "NAV Portugal (Continental)"
$ws2.Range("B30").Value = 276568
$ws2.Range("C30").Value = 105979
$ws2.Range("D30").Value = 152942
$ws2.Range("A31").Value = "NAVIAIR"
$ws2.Range("B31").Value = 178196
$ws2.Range("C31").Value = 94628
$ws2.Range("D31").Value = 121307
$ws2.Range("A32").Value = "Oro navigacija"
$ws2.Range("B32").Value = 74537
$ws2.Range("C32").Value = 45941
$ws2.Range("D32").Value = 54690
$ws2.Range("A33").Value = "PANSA"
$ws2.Range("B33").Value = 73375
$ws2.Range("C33").Value = 39198
$ws2.Range("D33").Value = 50061
$ws2.Range("A34").Value = "ROMATSA"
$ws2.Range("B34").Value = 104383
$ws2.Range("C34").Value = 108547
$ws2.Range("D34").Value = 107060
$ws2.Range("A35").Value = "Sakaeronavigatsia"
$ws2.Range("B35").Value = 25549
$ws2.Range("C35").Value = 14059
$ws2.Range("D35").Value = 15626
$ws2.Range("A36").Value = "Skeyes"
$ws2.Range("B36").Value = 211550
$ws2.Range("C36").Value = 149456
$ws2.Range("D36").Value = 164922
$ws2.Range("A37").Value = "Skyguide"
$ws2.Range("B37").Value = 228225
$ws2.Range("C37").Value = 173362
$ws2.Range("D37").Value = 184996
$ws2.Range("A38").Value = "Slovenia Control"
$ws2.Range("B38").Value = 113759
$ws2.Range("C38").Value = 80276
$ws2.Range("D38").Value = 93083
$ws2.Range("A39").Value = "SMATSA"
$ws2.Range("B39").Value = 65527
$ws2.Range("C39").Value = 43370
$ws2.Range("D39").Value = 51731
$ws2.Range("A40").Value = "EUROCONTROL area -average values"
$ws2.Range("B40").Value = 149078
$ws2.Range("C40").Value = 78632
$ws2.Range("D40").Value = 101304

# Remove now-unused number format (integer "0") from data cells
$ws2.Range("B2:D40").ClearFormats()

# --- Selection / active sheet state ---
$ws2.Range("K9").Select()
$ws1.Activate()
$ws1.Range("I5").Select()
